$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-like numeric strings in column D keep their original formatting
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "34.007.97"
$ws.Range("E2").Value = "  -2.42%  "
$ws.Range("D3").Value = "1.775.69"
$ws.Range("E3").Value = "  -1.71%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "220.60"
$ws.Range("E5").Value = "  -2.78%  "
$ws.Range("D6").Value = "0.546"
$ws.Range("E6").Value = "  -2.39%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("D8").Value = "31.21"
$ws.Range("E8").Value = "  -6.42%  "
$ws.Range("D9").Value = "0.285"
$ws.Range("E9").Value = "  -1.03%  "
$ws.Range("D10").Value = "0.0704"
$ws.Range("E10").Value = "  +4.35%  "
$ws.Range("D11").Value = "0.0921"
$ws.Range("E11").Value = "  -1.68%  "
$ws.Range("D12").Value = "2.031.77"
$ws.Range("E12").Value = "  -1.54%  "
$ws.Range("D13").Value = "1.770.05"
$ws.Range("E13").Value = "  -1.23%  "
$ws.Range("D14").Value = "10.52"
$ws.Range("E14").Value = "  -9.08%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "0.620"
$ws.Range("E15").Value = "  -3.92%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "33.905.72"
$ws.Range("E16").Value = "  -2.51%  "
$ws.Range("D17").Value = "4.19"
$ws.Range("E17").Value = "  -2.95%  "
$ws.Range("D18").Value = "67.73"
$ws.Range("E18").Value = "  -3.09%  "
$ws.Range("D19").Value = "243.09"
$ws.Range("E19").Value = "  -5.76%  "
$ws.Range("D20").Value = "0.0₃0772"
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("D21").Value = "1.00"
$ws.Range("E21").Value = "  +0.19%  "
$ws.Range("D22").Value = "10.50"
$ws.Range("E22").Value = "  -0.54%  "
$ws.Range("E23").Value = "  -5.51%  "
$ws.Range("E24").Value = "  -1.60%  "
$ws.Range("D25").Value = "158.03"
$ws.Range("E25").Value = "  -0.63%  "
$ws.Range("D26").Value = "16.32"
$ws.Range("E26").Value = "  -1.74%  "
$ws.Range("D27").Value = "6.98"
$ws.Range("E27").Value = "  -3.07%  "
$ws.Range("E28").Value = "  -3.07%  "
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  +0.25%  "
$ws.Range("D30").Value = "0.0515"
$ws.Range("E30").Value = "  -1.91%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "3.68"
$ws.Range("E31").Value = "  -4.15%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "1.19"
$ws.Range("E32").Value = "  -0.97%  "
$ws.Range("D33").Value = "3.50"
$ws.Range("E33").Value = "  -3.79%  "
$ws.Range("D34").Value = "1.83"
$ws.Range("E34").Value = "  -5.12%  "
$ws.Range("D35").Value = "1.395.23"
$ws.Range("E35").Value = "  -5.01%  "
$ws.Range("E36").Value = "  -1.62%  "
$ws.Range("D37").Value = "0.625"
$ws.Range("E37").Value = "  -2.44%  "
$ws.Range("E38").Value = "  -3.49%  "
$ws.Range("D39").Value = "0.930"
$ws.Range("E39").Value = "  +2.10%  "
$ws.Range("E40").Value = "  +0.09%  "
$ws.Range("D41").Value = "78.80"
$ws.Range("E41").Value = "  -6.52%  "
$ws.Range("E42").Value = "  -5.79%  "
$ws.Range("E43").Value = "  -0.95%  "
$ws.Range("D44").Value = "0.0490"
$ws.Range("E44").Value = "  -3.49%  "
$ws.Range("D45").Value = "5.82"
$ws.Range("E45").Value = "  -3.57%  "
$ws.Range("E46").Value = "  -0.90%  "
$ws.Range("D47").Value = "1.924.14"
$ws.Range("E47").Value = "  -1.99%  "
$ws.Range("D48").Value = "103.32"
$ws.Range("E48").Value = "  -0.02%  "
$ws.Range("D49").Value = "0.994"
$ws.Range("E49").Value = "  -0.35%  "
$ws.Range("D50").Value = "11.72"
$ws.Range("E50").Value = "  -3.28%  "
$ws.Range("D51").Value = "0.0₆0115"
$ws.Range("E51").Value = "  -5.18%  "
